$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final table values for rows 2-10 (header row 1 is unchanged).
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the numeric expression/specificity metrics.
# This reflects the recomputed NATMI LR-pairs table after adding the "ECs"
# sending/target cluster to the analysis (per commit: "Natmi following Dr Hou advice").
$data = @(
    @("ECs",  "Myoc", "Fzd3", "ECs",  2, 0.6666666666666666, 0.349434,            1.048302,            0.0009963999680650763, 0.0009963999680650763, 3, 1, 0.265744,           0.7972319999999999, 0.1034864391735229, 0.1034864391735229, 0.092859988896,     0.835739900064,     0.0001031138846876666, 0.0001031138846876666),
    @("ECs",  "Myoc", "Fzd3", "FAPs", 2, 0.6666666666666666, 0.349434,            1.048302,            0.0009963999680650763, 0.0009963999680650763, 3, 1, 0.307583,           0.922749,            0.1197794472135201, 0.1197794472135201, 0.107479958022,     0.9673196221980001, 0.0001193482373784039, 0.0001193482373784039),
    @("ECs",  "Myoc", "Fzd3", "sCs",  2, 0.6666666666666666, 0.349434,            1.048302,            0.0009963999680650763, 0.0009963999680650763, 3, 1, 1.994584333333333, 5.983753,            0.7767341136129571, 0.7767341136129571, 0.696975581934,     6.272780237406001,  0.0007739378459990057, 0.0007739378459990057),
    @("FAPs", "Myoc", "Fzd3", "ECs",  3, 1,                   348.977468,         1046.932404,         0.9950981815468188,    0.9950981815468188,    3, 1, 0.265744,           0.7972319999999999, 0.1034864391735229, 0.1034864391735229, 92.73866825619199,  834.6480143057279,  0.1029791674363281,    0.1029791674363281),
    @("FAPs", "Myoc", "Fzd3", "FAPs", 3, 1,                   348.977468,         1046.932404,         0.9950981815468188,    0.9950981815468188,    3, 1, 0.307583,           0.922749,            0.1197794472135201, 0.1197794472135201, 107.339536539844,   966.055828858596,   0.119192310108857,     0.119192310108857),
    @("FAPs", "Myoc", "Fzd3", "sCs",  3, 1,                   348.977468,         1046.932404,         0.9950981815468188,    0.9950981815468188,    3, 1, 1.994584333333333, 5.983753,            0.7767341136129571, 0.7767341136129571, 696.0649903591346,  6264.584913232212,  0.7729267040016338,    0.7729267040016338),
    @("sCs",  "Myoc", "Fzd3", "ECs",  3, 1,                   1.369616666666667,  4.10885,             0.00390541848511611,   0.00390541848511611,   3, 1, 0.265744,           0.7972319999999999, 0.1034864391735229, 0.1034864391735229, 0.3639674114666667, 3.2757067032,       0.0004041578525071201, 0.0004041578525071201),
    @("sCs",  "Myoc", "Fzd3", "FAPs", 3, 1,                   1.369616666666667,  4.10885,             0.00390541848511611,   0.00390541848511611,   3, 1, 0.307583,           0.922749,            0.1197794472135201, 0.1197794472135201, 0.4212708031833333, 3.79143722865,      0.0004677888672846707, 0.0004677888672846707),
    @("sCs",  "Myoc", "Fzd3", "sCs",  3, 1,                   1.369616666666667,  4.10885,             0.00390541848511611,   0.00390541848511611,   3, 1, 1.994584333333333, 5.983753,            0.7767341136129571, 0.7767341136129571, 2.731815946005556,  24.58634351405,     0.003033471765324319,  0.003033471765324319)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowValues[$j]
    }
}
